$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = '57.781.65'
$ws.Cells.Item(2, 5).Value = '  -4.62%  '

$ws.Cells.Item(3, 4).Value = '3.167.58'
$ws.Cells.Item(3, 5).Value = '  -5.30%  '

$ws.Cells.Item(4, 5).Value = '  -0.02%  '

$ws.Cells.Item(5, 4).Value = '529.37'
$ws.Cells.Item(5, 5).Value = '  -6.60%  '

$ws.Cells.Item(6, 4).Value = '135.27'
$ws.Cells.Item(6, 5).Value = '  -7.64%  '

$ws.Cells.Item(7, 5).Value = '  -0.09%  '

$ws.Cells.Item(8, 4).Value = '3.165.31'
$ws.Cells.Item(8, 5).Value = '  -5.41%  '

$ws.Cells.Item(9, 5).Value = '  -6.56%  '

$ws.Cells.Item(10, 4).Value = '7.22'
$ws.Cells.Item(10, 5).Value = '  -8.67%  '

$ws.Cells.Item(11, 5).Value = '  -8.16%  '

$ws.Cells.Item(12, 4).Value = '0.395'
$ws.Cells.Item(12, 5).Value = '  -4.80%  '

$ws.Cells.Item(13, 4).Value = '3.707.37'
$ws.Cells.Item(13, 5).Value = '  -5.31%  '

$ws.Cells.Item(14, 5).Value = '  -1.49%  '

$ws.Cells.Item(15, 4).Value = '25.94'
$ws.Cells.Item(15, 5).Value = '  -6.30%  '

$ws.Cells.Item(16, 4).Value = '3.161.27'
$ws.Cells.Item(16, 5).Value = '  -5.19%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(17, 4).Value = '57.710.76'
$ws.Cells.Item(17, 5).Value = '  -4.78%  '

$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.0000154'
$ws.Cells.Item(18, 5).Value = '  -8.73%  '

$ws.Cells.Item(19, 5).Value = '  -6.95%  '

$ws.Cells.Item(20, 4).Value = '13.17'
$ws.Cells.Item(20, 5).Value = '  -9.43%  '

$ws.Cells.Item(21, 4).Value = '8.10'
$ws.Cells.Item(21, 5).Value = '  -8.86%  '

$ws.Cells.Item(22, 4).Value = '349.92'
$ws.Cells.Item(22, 5).Value = '  -7.04%  '

$ws.Cells.Item(23, 5).Value = '  -0.02%  '

$ws.Cells.Item(24, 4).Value = '69.69'
$ws.Cells.Item(24, 5).Value = '  -6.74%  '

$ws.Cells.Item(25, 4).Value = '0.514'
$ws.Cells.Item(25, 5).Value = '  -8.03%  '

$ws.Cells.Item(26, 4).Value = '3.290.42'
$ws.Cells.Item(26, 5).Value = '  -5.72%  '

$ws.Cells.Item(27, 4).Value = '0.0₃0970'
$ws.Cells.Item(27, 5).Value = '  -10.13%  '

$ws.Cells.Item(28, 5).Value = '  -3.75%  '

$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  -0.25%  '

$ws.Cells.Item(30, 4).Value = '6.96'
$ws.Cells.Item(30, 5).Value = '  -5.14%  '

$ws.Cells.Item(31, 5).Value = '  -0.11%  '

$ws.Cells.Item(32, 5).Value = '  -9.17%  '

$ws.Cells.Item(33, 4).Value = '6.99'
$ws.Cells.Item(33, 5).Value = '  -8.93%  '

$ws.Cells.Item(34, 4).Value = '21.72'
$ws.Cells.Item(34, 5).Value = '  -5.09%  '

$ws.Cells.Item(35, 5).Value = '  -5.53%  '

$ws.Cells.Item(36, 4).Value = '4.96'
$ws.Cells.Item(36, 5).Value = '  -6.18%  '

$ws.Cells.Item(37, 4).Value = '158.84'
$ws.Cells.Item(37, 5).Value = '  -5.08%  '

$ws.Cells.Item(38, 2).Value = 'Aptos'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(38, 4).Value = '6.28'
$ws.Cells.Item(38, 5).Value = '  -7.75%  '

$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = '1.42'
$ws.Cells.Item(39, 5).Value = '  -8.64%  '

$ws.Cells.Item(40, 4).Value = '26.29'
$ws.Cells.Item(40, 5).Value = '  -5.92%  '

$ws.Cells.Item(41, 4).Value = '0.0703'
$ws.Cells.Item(41, 5).Value = '  -5.79%  '

$ws.Cells.Item(42, 4).Value = '3.194.29'
$ws.Cells.Item(42, 5).Value = '  -5.50%  '

$ws.Cells.Item(43, 4).Value = '40.37'
$ws.Cells.Item(43, 5).Value = '  -4.28%  '

$ws.Cells.Item(44, 2).Value = 'Mantle'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(44, 4).Value = '0.699'
$ws.Cells.Item(44, 5).Value = '  -7.71%  '

$ws.Cells.Item(45, 2).Value = 'ONDO'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(45, 4).Value = '1.10'
$ws.Cells.Item(45, 5).Value = '  -3.90%  '

$ws.Cells.Item(46, 4).Value = '3.99'
$ws.Cells.Item(46, 5).Value = '  -6.88%  '

$ws.Cells.Item(47, 4).Value = '1.00'
$ws.Cells.Item(47, 5).Value = '  -0.12%  '

$ws.Cells.Item(48, 5).Value = '  -8.35%  '

$ws.Cells.Item(49, 4).Value = '2.276.98'
$ws.Cells.Item(49, 5).Value = '  -7.22%  '

$ws.Cells.Item(50, 5).Value = '  -6.55%  '

$ws.Cells.Item(51, 4).Value = '20.79'
$ws.Cells.Item(51, 5).Value = '  -6.91%  '
